$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("massaDados")

# --- Row 2: JOÃO/SILVA -> ESDSON/MARCANTE (Business Name stays the same) ---
$ws.Range("A2").Value = "ESDSON"
$ws.Range("B2").Value = "MARCANTE"

# --- Rows 4-11: clear the data (name/last name/business) and the email hyperlink,
#     leaving just the formatted (blank) D column cell behind ---
$ws.Range("A4:C11").ClearContents()
$ws.Range("D4:D11").ClearContents()

# Hyperlinks.Delete() on this host clears the *entire* sheet collection rather
# than just the addressed range, so remove them all and re-create only the two
# that must survive (D2, D3), then restore their original "Hiperlink" cell
# style (Add() otherwise stamps a freshly-applied-font variant of it).
$ws.Range("A1:E11").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:DOZE@GMAIL.COM")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:DOZE@GMAIL.COM")
$ws.Range("D2").Style = "Hiperlink"
$ws.Range("D3").Style = "Hiperlink"

# --- Selection moves to A4:E11 with A4 active ---
$ws.Range("A4:E11").Select()
